# Applies odds updates to Sheet1 for 2025-12-09 Betfair Back/Lay data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("AC3").Value = 9.6
$ws.Range("AD3").Value = 25
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 13
$ws.Range("AH3").Value = 27
$ws.Range("AJ3").Value = 27
$ws.Range("AK3").Value = 28
$ws.Range("AL3").Value = 60
$ws.Range("AN3").Value = 21
$ws.Range("F3").Value = 1.8
$ws.Range("G3").Value = 2.02
$ws.Range("H3").Value = 4.7
$ws.Range("K3").Value = 3.9
$ws.Range("L3").Value = 1.46
$ws.Range("N3").Value = 2.96
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 1.67
$ws.Range("Q3").Value = 2.16
$ws.Range("R3").Value = 1.25
$ws.Range("S3").Value = 4.1
$ws.Range("T3").Value = 1.97
$ws.Range("U3").Value = 1.81
$ws.Range("V3").Value = 1.18
$ws.Range("W3").Value = 1.98
$ws.Range("X3").Value = 13.5
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 48

# Row 4
$ws.Range("N4").Value = 1.1
$ws.Range("P4").Value = 1.52
$ws.Range("R4").Value = 1.17

# Row 5
$ws.Range("AA5").Value = 13
$ws.Range("F5").Value = 7.4
$ws.Range("G5").Value = 7.6
$ws.Range("Q5").Value = 1.78
$ws.Range("R5").Value = 1.47

# Row 6
$ws.Range("AG6").Value = 11
$ws.Range("N6").Value = 7.8

# Row 7
$ws.Range("O7").Value = 1.24
$ws.Range("S7").Value = 2.86
$ws.Range("U7").Value = 2.46

# Row 8
$ws.Range("AH8").Value = 27
$ws.Range("AI8").Value = 1000
$ws.Range("AO8").Value = 180
$ws.Range("I8").Value = 8.6
$ws.Range("J8").Value = 5.1

# Row 9
$ws.Range("AC9").Value = 9.2
$ws.Range("AH9").Value = 14
$ws.Range("AN9").Value = 21
$ws.Range("M9").Value = 1.04
$ws.Range("P9").Value = 2.56

# Row 10
$ws.Range("Y10").Value = 19.5

# Row 11
$ws.Range("AC11").Value = 8.6
$ws.Range("G11").Value = 2.22
$ws.Range("T11").Value = 1.59
$ws.Range("W11").Value = 1.81

# Row 12
$ws.Range("AF12").Value = 12
$ws.Range("AK12").Value = 13
$ws.Range("H12").Value = 17.5
$ws.Range("P12").Value = 4.3
$ws.Range("T12").Value = 1.82
